$wb = $excel.ActiveWorkbook

# --- Add 4 new worksheets after Physical_contacts_imputed ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHome = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsHome.Name = "All_contacts_home"

$wsWork = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsHome)
$wsWork.Name = "All_contacts_work"

$wsSchool = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsWork)
$wsSchool.Name = "All_contacts_school"

$wsOther = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSchool)
$wsOther.Name = "All_contacts_other"

# --- Populate All_contacts_home ---
$arr6 = New-Object 'object[,]' 9,9
$arr6[0,1] = "[0,5)"
$arr6[0,2] = "[5,18)"
$arr6[0,3] = "[18,30)"
$arr6[0,4] = "[30,40)"
$arr6[0,5] = "[40,50)"
$arr6[0,6] = "[50,60)"
$arr6[0,7] = "[60,70)"
$arr6[0,8] = "70+"
$arr6[1,0] = "[0,5)"
$arr6[1,1] = "NA"
$arr6[1,2] = "NA"
$arr6[1,3] = "NA"
$arr6[1,4] = "NA"
$arr6[1,5] = "NA"
$arr6[1,6] = "NA"
$arr6[1,7] = "NA"
$arr6[1,8] = "NA"
$arr6[2,0] = "[5,18)"
$arr6[2,1] = "NA"
$arr6[2,2] = "NA"
$arr6[2,3] = "NA"
$arr6[2,4] = "NA"
$arr6[2,5] = "NA"
$arr6[2,6] = "NA"
$arr6[2,7] = "NA"
$arr6[2,8] = "NA"
$arr6[3,0] = "[18,30)"
$arr6[3,1] = 0.097560975609756101
$arr6[3,2] = 0.27804878048780501
$arr6[3,3] = 0.84878048780487803
$arr6[3,4] = 0.063414634146341506
$arr6[3,5] = 0.30731707317073198
$arr6[3,6] = 0.27317073170731698
$arr6[3,7] = 0.0292682926829268
$arr6[3,8] = 0.034146341463414602
$arr6[4,0] = "[30,40)"
$arr6[4,1] = 0.178988326848249
$arr6[4,2] = 0.39688715953307402
$arr6[4,3] = 0.57976653696498104
$arr6[4,4] = 0.369649805447471
$arr6[4,5] = 0.062256809338521402
$arr6[4,6] = 0.101167315175097
$arr6[4,7] = 0.050583657587548597
$arr6[4,8] = 0.054474708171206199
$arr6[5,0] = "[40,50)"
$arr6[5,1] = 0.079295154185022004
$arr6[5,2] = 0.64317180616740099
$arr6[5,3] = 0.14977973568281899
$arr6[5,4] = 0.3215859030837
$arr6[5,5] = 0.45374449339207001
$arr6[5,6] = 0.070484581497797405
$arr6[5,7] = 0.035242290748898703
$arr6[5,8] = 0.052863436123347998
$arr6[6,0] = "[50,60)"
$arr6[6,1] = 0.013824884792626699
$arr6[6,2] = 0.29493087557603698
$arr6[6,3] = 0.26728110599078297
$arr6[6,4] = 0.082949308755760398
$arr6[6,5] = 0.43778801843317999
$arr6[6,6] = 0.33179723502304098
$arr6[6,7] = 0.032258064516128997
$arr6[6,8] = 0.101382488479263
$arr6[7,0] = "[60,70)"
$arr6[7,1] = 0.00377358490566038
$arr6[7,2] = 0.045283018867924497
$arr6[7,3] = 0.13962264150943399
$arr6[7,4] = 0.109433962264151
$arr6[7,5] = 0.045283018867924497
$arr6[7,6] = 0.43018867924528298
$arr6[7,7] = 0.354716981132075
$arr6[7,8] = 0.128301886792453
$arr6[8,0] = "70+"
$arr6[8,1] = 0
$arr6[8,2] = 0.020833333333333301
$arr6[8,3] = 0.010416666666666701
$arr6[8,4] = 0.03125
$arr6[8,5] = 0.072916666666666699
$arr6[8,6] = 0.0625
$arr6[8,7] = 0.14583333333333301
$arr6[8,8] = 0.77083333333333304
$wsHome.Range("A1:I9").Value2 = $arr6
$wsHome.Range("M33").Select()

# --- Populate All_contacts_work ---
$arr7 = New-Object 'object[,]' 9,9
$arr7[0,1] = "[0,5)"
$arr7[0,2] = "[5,18)"
$arr7[0,3] = "[18,30)"
$arr7[0,4] = "[30,40)"
$arr7[0,5] = "[40,50)"
$arr7[0,6] = "[50,60)"
$arr7[0,7] = "[60,70)"
$arr7[0,8] = "70+"
$arr7[1,0] = "[0,5)"
$arr7[1,1] = "NA"
$arr7[1,2] = "NA"
$arr7[1,3] = "NA"
$arr7[1,4] = "NA"
$arr7[1,5] = "NA"
$arr7[1,6] = "NA"
$arr7[1,7] = "NA"
$arr7[1,8] = "NA"
$arr7[2,0] = "[5,18)"
$arr7[2,1] = "NA"
$arr7[2,2] = "NA"
$arr7[2,3] = "NA"
$arr7[2,4] = "NA"
$arr7[2,5] = "NA"
$arr7[2,6] = "NA"
$arr7[2,7] = "NA"
$arr7[2,8] = "NA"
$arr7[3,0] = "[18,30)"
$arr7[3,1] = 0.0096618357487922701
$arr7[3,2] = 0.043478260869565202
$arr7[3,3] = 0.15942028985507201
$arr7[3,4] = 0.101449275362319
$arr7[3,5] = 0.038647342995169101
$arr7[3,6] = 0.028985507246376802
$arr7[3,7] = 0.019323671497584499
$arr7[3,8] = 0.0096618357487922701
$arr7[4,0] = "[30,40)"
$arr7[4,1] = 0.0037037037037036999
$arr7[4,2] = 0.033333333333333298
$arr7[4,3] = 0.23703703703703699
$arr7[4,4] = 0.203703703703704
$arr7[4,5] = 0.16666666666666699
$arr7[4,6] = 0.077777777777777807
$arr7[4,7] = 0.0148148148148148
$arr7[4,8] = 0.0148148148148148
$arr7[5,0] = "[40,50)"
$arr7[5,1] = 0
$arr7[5,2] = 0.0041493775933610002
$arr7[5,3] = 0.128630705394191
$arr7[5,4] = 0.14937759336099601
$arr7[5,5] = 0.182572614107884
$arr7[5,6] = 0.141078838174274
$arr7[5,7] = 0.00829875518672199
$arr7[5,8] = 0
$arr7[6,0] = "[50,60)"
$arr7[6,1] = 0
$arr7[6,2] = 0.0043859649122806998
$arr7[6,3] = 0.118421052631579
$arr7[6,4] = 0.16666666666666699
$arr7[6,5] = 0.162280701754386
$arr7[6,6] = 0.157894736842105
$arr7[6,7] = 0.0087719298245613996
$arr7[6,8] = 0.021929824561403501
$arr7[7,0] = "[60,70)"
$arr7[7,1] = 0
$arr7[7,2] = 0
$arr7[7,3] = 0.054151624548736503
$arr7[7,4] = 0.086642599277978294
$arr7[7,5] = 0.064981949458483707
$arr7[7,6] = 0.046931407942238303
$arr7[7,7] = 0.0072202166064982004
$arr7[7,8] = 0.014440433212996401
$arr7[8,0] = "70+"
$arr7[8,1] = 0
$arr7[8,2] = 0
$arr7[8,3] = 0
$arr7[8,4] = 0
$arr7[8,5] = 0
$arr7[8,6] = 0
$arr7[8,7] = 0
$arr7[8,8] = 0
$wsWork.Range("A1:I9").Value2 = $arr7
$wsWork.Range("A1:I9").Select()

# --- Populate All_contacts_school ---
$arr8 = New-Object 'object[,]' 9,9
$arr8[0,1] = "[0,5)"
$arr8[0,2] = "[5,18)"
$arr8[0,3] = "[18,30)"
$arr8[0,4] = "[30,40)"
$arr8[0,5] = "[40,50)"
$arr8[0,6] = "[50,60)"
$arr8[0,7] = "[60,70)"
$arr8[0,8] = "70+"
$arr8[1,0] = "[0,5)"
$arr8[1,1] = "NA"
$arr8[1,2] = "NA"
$arr8[1,3] = "NA"
$arr8[1,4] = "NA"
$arr8[1,5] = "NA"
$arr8[1,6] = "NA"
$arr8[1,7] = "NA"
$arr8[1,8] = "NA"
$arr8[2,0] = "[5,18)"
$arr8[2,1] = "NA"
$arr8[2,2] = "NA"
$arr8[2,3] = "NA"
$arr8[2,4] = "NA"
$arr8[2,5] = "NA"
$arr8[2,6] = "NA"
$arr8[2,7] = "NA"
$arr8[2,8] = "NA"
$arr8[3,0] = "[18,30)"
$arr8[3,1] = 0.0046948356807511703
$arr8[3,2] = 0.0093896713615023494
$arr8[3,3] = 0.0093896713615023494
$arr8[3,4] = 0.0093896713615023494
$arr8[3,5] = 0
$arr8[3,6] = 0.0046948356807511703
$arr8[3,7] = 0.0046948356807511703
$arr8[3,8] = 0
$arr8[4,0] = "[30,40)"
$arr8[4,1] = 0.0036363636363636398
$arr8[4,2] = 0.0036363636363636398
$arr8[4,3] = 0
$arr8[4,4] = 0
$arr8[4,5] = 0
$arr8[4,6] = 0.0036363636363636398
$arr8[4,7] = 0
$arr8[4,8] = 0.0072727272727272701
$arr8[5,0] = "[40,50)"
$arr8[5,1] = 0
$arr8[5,2] = 0
$arr8[5,3] = 0
$arr8[5,4] = 0
$arr8[5,5] = 0
$arr8[5,6] = 0
$arr8[5,7] = 0
$arr8[5,8] = 0
$arr8[6,0] = "[50,60)"
$arr8[6,1] = 0
$arr8[6,2] = 0
$arr8[6,3] = 0
$arr8[6,4] = 0
$arr8[6,5] = 0
$arr8[6,6] = 0
$arr8[6,7] = 0
$arr8[6,8] = 0
$arr8[7,0] = "[60,70)"
$arr8[7,1] = 0
$arr8[7,2] = 0
$arr8[7,3] = 0
$arr8[7,4] = 0
$arr8[7,5] = 0
$arr8[7,6] = 0
$arr8[7,7] = 0
$arr8[7,8] = 0
$arr8[8,0] = "70+"
$arr8[8,1] = 0
$arr8[8,2] = 0
$arr8[8,3] = 0
$arr8[8,4] = 0
$arr8[8,5] = 0
$arr8[8,6] = 0
$arr8[8,7] = 0
$arr8[8,8] = 0
$wsSchool.Range("A1:I9").Value2 = $arr8
$wsSchool.Range("A1:I9").Select()

# --- Populate All_contacts_other ---
$arr9 = New-Object 'object[,]' 9,9
$arr9[0,1] = "[0,5)"
$arr9[0,2] = "[5,18)"
$arr9[0,3] = "[18,30)"
$arr9[0,4] = "[30,40)"
$arr9[0,5] = "[40,50)"
$arr9[0,6] = "[50,60)"
$arr9[0,7] = "[60,70)"
$arr9[0,8] = "70+"
$arr9[1,0] = "[0,5)"
$arr9[1,1] = "NA"
$arr9[1,2] = "NA"
$arr9[1,3] = "NA"
$arr9[1,4] = "NA"
$arr9[1,5] = "NA"
$arr9[1,6] = "NA"
$arr9[1,7] = "NA"
$arr9[1,8] = "NA"
$arr9[2,0] = "[5,18)"
$arr9[2,1] = "NA"
$arr9[2,2] = "NA"
$arr9[2,3] = "NA"
$arr9[2,4] = "NA"
$arr9[2,5] = "NA"
$arr9[2,6] = "NA"
$arr9[2,7] = "NA"
$arr9[2,8] = "NA"
$arr9[3,0] = "[18,30)"
$arr9[3,1] = 0.0153061224489796
$arr9[3,2] = 0.025510204081632699
$arr9[3,3] = 0.17857142857142899
$arr9[3,4] = 0.035714285714285698
$arr9[3,5] = 0.076530612244898003
$arr9[3,6] = 0.025510204081632699
$arr9[3,7] = 0.0102040816326531
$arr9[3,8] = 0.0204081632653061
$arr9[4,0] = "[30,40)"
$arr9[4,1] = 0.020746887966804999
$arr9[4,2] = 0.024896265560166001
$arr9[4,3] = 0.10788381742738599
$arr9[4,4] = 0.058091286307053902
$arr9[4,5] = 0.0539419087136929
$arr9[4,6] = 0.045643153526971
$arr9[4,7] = 0.00829875518672199
$arr9[4,8] = 0.029045643153527
$arr9[5,0] = "[40,50)"
$arr9[5,1] = 0.013698630136986301
$arr9[5,2] = 0.022831050228310501
$arr9[5,3] = 0.082191780821917804
$arr9[5,4] = 0.077625570776255703
$arr9[5,5] = 0.11872146118721499
$arr9[5,6] = 0.054794520547945202
$arr9[5,7] = 0.027397260273972601
$arr9[5,8] = 0.041095890410958902
$arr9[6,0] = "[50,60)"
$arr9[6,1] = 0.0104712041884817
$arr9[6,2] = 0.026178010471204199
$arr9[6,3] = 0.104712041884817
$arr9[6,4] = 0.083769633507853394
$arr9[6,5] = 0.089005235602094196
$arr9[6,6] = 0.094240837696335095
$arr9[6,7] = 0.0104712041884817
$arr9[6,8] = 0.083769633507853394
$arr9[7,0] = "[60,70)"
$arr9[7,1] = 0.013333333333333299
$arr9[7,2] = 0
$arr9[7,3] = 0.124444444444444
$arr9[7,4] = 0.084444444444444405
$arr9[7,5] = 0.10222222222222201
$arr9[7,6] = 0.11555555555555599
$arr9[7,7] = 0.093333333333333296
$arr9[7,8] = 0.053333333333333302
$arr9[8,0] = "70+"
$arr9[8,1] = 0.047619047619047603
$arr9[8,2] = 0
$arr9[8,3] = 0.107142857142857
$arr9[8,4] = 0.095238095238095205
$arr9[8,5] = 0.0595238095238095
$arr9[8,6] = 0.071428571428571397
$arr9[8,7] = 0.083333333333333301
$arr9[8,8] = 0.214285714285714
$wsOther.Range("A1:I9").Value2 = $arr9
$wsOther.Range("A1:I9").Select()

# --- Update view/selection on existing sheets ---
# Physical_contacts: was the selected tab; selection itself is unchanged,
# but it loses tabSelected once Detail is activated below.
$wsPhysical = $wb.Worksheets.Item("Physical_contacts")
$wsPhysical.Range("B12:H19").Select()

# Detail becomes the selected tab, with the active cell moved to A10.
$wsDetail = $wb.Worksheets.Item("Detail")
$wsDetail.Range("A10").Select()
$wsDetail.Activate()
